# Add a new "Datetime Components" worksheet after the "Regex (Regular Expressions)"
# sheet, matching the structure used by the other "GOMS Python results" task sheets
# (Task / Action-Time-Content table, rows 1..6).

$wb = $excel.ActiveWorkbook

$regexSheet = $wb.Worksheets.Item("Regex (Regular Expressions)")

# Insert the new sheet right after the Regex sheet (it becomes the new last sheet
# and the active tab, just like the Regex sheet was before this edit).
$newSheet = $wb.Worksheets.Add($null, $regexSheet)
$newSheet.Name = "Datetime Components"

# --- Row 1: header (bold) ---------------------------------------------------
$newSheet.Range("A1").Value = "Action"
$newSheet.Range("B1").Value = "Time"
$newSheet.Range("C1").Value = "Content"
$newSheet.Range("A1:C1").Font.Size = 13
$newSheet.Range("A1:C1").Font.Bold = $true

# --- Row 2: Upload CSV -------------------------------------------------------
$newSheet.Range("A2").Value = "Upload CSV"
$newSheet.Range("B2").Value = "5 min"
$newSheet.Range("C2").Value = "df = pd.read_csv('file.csv')"
$newSheet.Range("A2:C2").Font.Size = 13

# --- Row 3: Convert to Datetime ---------------------------------------------
$newSheet.Range("A3").Value = "Convert to Datetime"
$newSheet.Range("B3").Value = "2 min"
$newSheet.Range("C3").Value = "df['datetime_column'] = pd.to_datetime(df['datetime_column'])"
$newSheet.Range("A3:C3").Font.Size = 13

# --- Row 4: Extract Components ----------------------------------------------
$newSheet.Range("A4").Value = "Extract Components"
$newSheet.Range("B4").Value = "3 min"
$newSheet.Range("C4").Value = "df['year'] = df['datetime_column'].dt.year etc. for month, day, etc."
$newSheet.Range("A4:C4").Font.Size = 13

# --- Row 5: Verify Changes ---------------------------------------------------
$newSheet.Range("A5").Value = "Verify Changes"
$newSheet.Range("B5").Value = "1 min"
$newSheet.Range("C5").Value = "df[['year', 'month', 'day']].head()"
$newSheet.Range("A5:C5").Font.Size = 13

# --- Row 6: Overall (bold) ---------------------------------------------------
$newSheet.Range("A6").Value = "Overall"
$newSheet.Range("B6").Value = "11 min"
$newSheet.Range("A6:B6").Font.Size = 13
$newSheet.Range("A6:B6").Font.Bold = $true

# Row heights to match the other task sheets (17pt rows).
$newSheet.Rows.Item(1).RowHeight = 17
$newSheet.Rows.Item(2).RowHeight = 17
$newSheet.Rows.Item(3).RowHeight = 17
$newSheet.Rows.Item(4).RowHeight = 17
$newSheet.Rows.Item(5).RowHeight = 17
$newSheet.Rows.Item(6).RowHeight = 17

# Match selection/appearance used on the other task sheets.
$newSheet.Range("A1:C6").Select() | Out-Null
